# Commit: "the babies with 2 or more files now have only one"
#
# For each baby that currently has two rows (one per recording file,
# labelled "<id>_1" / "<id>_2" in column A), insert a new combined row
# right after the pair that merges the two files into one: most columns
# become the AVERAGE of the pair, column J (number of spikes) becomes the
# SUM of the pair, and column A gets the baby's plain numeric id.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (other than A and J) that get AVERAGE(<col><r1>,<col><r2>)
$avgCols = @("B","C","D","E","F","G","H","I","K","L","M","N","O")

# Pairs of rows (as they exist in the sheet *before* any insertion in this
# loop), each being the "<id>_1" row followed immediately by the "<id>_2"
# row. Processed top-to-bottom, adjusting for rows already inserted above.
$pairs = @(
    @{R1 = 16; R2 = 17},
    @{R1 = 26; R2 = 27},
    @{R1 = 35; R2 = 36},
    @{R1 = 49; R2 = 50}
)

$shift = 0
foreach ($pair in $pairs) {
    $r1 = $pair.R1 + $shift
    $r2 = $pair.R2 + $shift
    $insertAt = $r2 + 1

    $ws.Rows.Item($insertAt).Insert()

    foreach ($col in $avgCols) {
        $ws.Range("$col$insertAt").Formula = "=AVERAGE($col$r1,$col$r2)"
    }
    $ws.Range("J$insertAt").Formula = "=J$r1+J$r2"

    # Column A: baby id without the "_1"/"_2" file suffix, as a number.
    $idText = [string]$ws.Range("A$r1").Value2
    $babyId = [int]($idText.Split("_")[0])
    $ws.Range("A$insertAt").Value = $babyId

    $shift = $shift + 1
}

# Update the visible selection to match the edited state.
$ws.Range("A54:O54").Select()
